$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 64.09999999999999
$ws.Range("K3").Value = 54.7
$ws.Range("K4").Value = 52.9
$ws.Range("K5").Value = 52.7
$ws.Range("K6").Value = 48.9

$ws.Range("N2").Value = 85.82376350509293
$ws.Range("N3").Value = 85.82376350509293
$ws.Range("N4").Value = 85.82376350509293
$ws.Range("N5").Value = 85.82376350509293
$ws.Range("N6").Value = 85.82376350509293
